# Update countries & provincias Spain
# Refresh of the COVID "Pais" data sheet: a handful of country rows got new
# totals (which re-sorted them against their neighbour), plus the
# "Datos actualizados" timestamp advanced by an hour.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp banner (A1) --------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 08:05"

# --- India (row 10) : new cases reported, overall total/recovered bumped ---
$ws.Range("B10").Value = 190791
$ws.Range("C10").Value = 182
$ws.Range("E10").Value = 93528

# --- Israel (row 44) : new cases reported -----------------------------------
$ws.Range("B44").Value = 17106
$ws.Range("C44").Value = 35
$ws.Range("D44").Value = 14826
$ws.Range("E44").Value = 1995

# --- Uzbekistan overtakes Senegal (rows 78/79) ------------------------------
$ws.Range("A78").Value = "Uzbekistan"
$ws.Range("B78").Value = 3662
$ws.Range("C78").Value = 39
$ws.Range("D78").Value = 2837
$ws.Range("E78").Value = 810
$ws.Range("H78").Value = 15

$ws.Range("A79").Value = "Senegal"
$ws.Range("B79").Value = 3645
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 1801
$ws.Range("E79").Value = 1802
$ws.Range("H79").Value = 42

# --- Bulgaria overtakes El Salvador (rows 86/87) ----------------------------
$ws.Range("A86").Value = "Bulgaria"
$ws.Range("B86").Value = 2519
$ws.Range("C86").Value = 6
$ws.Range("D86").Value = 1090
$ws.Range("E86").Value = 1289
$ws.Range("H86").Value = 140

$ws.Range("A87").Value = "El Salvador"
$ws.Range("B87").Value = 2517
$ws.Range("C87").Value = 0
$ws.Range("D87").Value = 1040
$ws.Range("E87").Value = 1431
$ws.Range("H87").Value = 46

# --- Belice overtakes Santa Lucia (rows 201/202) ----------------------------
$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2

$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("D202").Value = 18
$ws.Range("H202").Value = 0

# --- Seychelles overtakes Montserrat (rows 210/211) -------------------------
$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# --- San Bartolome / Bonaire tie swaps order (rows 215/216, values equal) --
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"
